# UI changes and method implementation
# Implemented stop, previous, next methods and some UI changes.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the grading values (G4, G5, G11, G12, G18, G19, G21).
# Dependent formulas (G6, H6, G13, H13, G22, H22, H32, I32, I34) recalc automatically.
$ws.Range("G4").Value = 2
$ws.Range("G5").Value = 12
$ws.Range("G11").Value = 3
$ws.Range("G12").Value = 6
$ws.Range("G18").Value = 1
$ws.Range("G19").Value = 2
$ws.Range("G21").Value = 2

# Move the active selection to reflect where the user last worked.
$ws.Range("G24").Select()
